$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header updates
$ws.Range("B2").Value = "16-07-24"
$ws.Range("B3").Value = "PARLE PRODUCTS"
$ws.Range("B4").Value = "470MM HM SHEET 4008115"

# Update data rows 6-9 (keep text formatting like the original cells)
$ws.Range("B6:F9").NumberFormat = "@"

$ws.Range("B6").Value = "30"
$ws.Range("C6").Value = "26.210"
$ws.Range("D6").Value = "0.08"
$ws.Range("E6").Value = "0.9"
$ws.Range("F6").Value = "25.230"

$ws.Range("B7").Value = "31"
$ws.Range("C7").Value = "27.480"
$ws.Range("D7").Value = "0.08"
$ws.Range("E7").Value = "0.9"
$ws.Range("F7").Value = "26.500"

$ws.Range("B8").Value = "32"
$ws.Range("C8").Value = "33.530"
$ws.Range("D8").Value = "0.08"
$ws.Range("E8").Value = "0.9"
$ws.Range("F8").Value = "32.550"

$ws.Range("B9").Value = "33"
$ws.Range("C9").Value = "26.360"
$ws.Range("D9").Value = "0.08"
$ws.Range("E9").Value = "0.9"
$ws.Range("F9").Value = "25.380"

# Remove the now-unneeded extra data rows (old rows 10-15); this shifts the
# old Total row (row 16, with its already-blank A cell) up to become row 10.
$ws.Range("A10:A15").EntireRow.Delete()

# Set new totals (A10 stays blank, as it already was on the old Total row)
$ws.Range("B10").Value = "Total"
$ws.Range("C10").Value = 113.58
$ws.Range("D10").Value = 0.32
$ws.Range("E10").Value = 3.6
$ws.Range("F10").Value = 109.66

$wb.Save()
